$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H51").Value = 8727
$ws.Range("I51").Value = 9090.5
$ws.Range("K51").Value = 9090.5
$ws.Range("M51").Value = -8606.5
$ws.Range("H69").Value = 13999.5
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 13999.5
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H86").Value = 5727.278
$ws.Range("J86").Value = 5742.6665
$ws.Range("L86").Value = 5742.6665
$ws.Range("N86").Value = -7988.6665
$ws.Range("H89").Value = 5727.278
$ws.Range("J89").Value = 5742.6665
$ws.Range("L89").Value = 28713.3325
$ws.Range("N89").Value = -39945.3325
$ws.Range("H103").Value = 634.5714
$ws.Range("J103").Value = 488.8
$ws.Range("L103").Value = 1466.4
$ws.Range("N103").Value = -2638.4
$ws.Range("H106").Value = 4984.8
$ws.Range("I106").Value = 5093.6665
$ws.Range("J106").Value = 4005
$ws.Range("K106").Value = 5093.6665
$ws.Range("L106").Value = 4005
$ws.Range("M106").Value = -4462.6665
$ws.Range("N106").Value = -5267
$ws.Range("H116").Value = 8250
$ws.Range("I116").Value = 8000
$ws.Range("J116").Value = 8500
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 8500
$ws.Range("M116").Value = -4558
$ws.Range("N116").Value = -15384
$ws.Range("H138").Value = 2155.56
$ws.Range("J138").Value = 2447.6956
$ws.Range("L138").Value = 7343.0868
$ws.Range("N138").Value = -17623.0868
$ws.Range("H141").Value = 4289.8335
$ws.Range("I141").Value = 5021.625
$ws.Range("K141").Value = 15064.875
$ws.Range("M141").Value = -9884.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2059.4546
$ws.Range("I2").Value = 2081.75
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 2081.75
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -1968.75
$ws.Range("N2").Value = -2226
$ws.Range("H61").Value = 63076.42
$ws.Range("I61").Value = 1437
$ws.Range("J61").Value = 196628.5
$ws.Range("K61").Value = 1437
$ws.Range("L61").Value = 196628.5
$ws.Range("M61").Value = -1225
$ws.Range("N61").Value = -197052.5
$ws.Range("H116").Value = 2059.4546
$ws.Range("I116").Value = 2081.75
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 2081.75
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 212.25
$ws.Range("N116").Value = -6588
$ws.Range("H136").Value = 63076.42
$ws.Range("I136").Value = 1437
$ws.Range("J136").Value = 196628.5
$ws.Range("K136").Value = 4311
$ws.Range("L136").Value = 589885.5
$ws.Range("M136").Value = -1761
$ws.Range("N136").Value = -594985.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2059.4546
$ws.Range("I3").Value = 2081.75
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 2081.75
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1967.75
$ws.Range("N3").Value = -2228
$ws.Range("H99").Value = 20919.041
$ws.Range("J99").Value = 6750
$ws.Range("L99").Value = 6750
$ws.Range("N99").Value = -9746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10515.489
$ws.Range("I58").Value = 4185.8184
$ws.Range("J58").Value = 25435.428
$ws.Range("K58").Value = 4185.8184
$ws.Range("L58").Value = 25435.428
$ws.Range("M58").Value = -3982.8184
$ws.Range("N58").Value = -25841.428
$ws.Range("H105").Value = 7528.533
$ws.Range("I105").Value = 9227.416999999999
$ws.Range("J105").Value = 733
$ws.Range("K105").Value = 9227.416999999999
$ws.Range("L105").Value = 733
$ws.Range("M105").Value = -7480.416999999999
$ws.Range("N105").Value = -4227
$ws.Range("H136").Value = 10515.489
$ws.Range("I136").Value = 4185.8184
$ws.Range("J136").Value = 25435.428
$ws.Range("K136").Value = 12557.4552
$ws.Range("L136").Value = 76306.284
$ws.Range("M136").Value = -10007.4552
$ws.Range("N136").Value = -81406.284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 358.14285
$ws.Range("I23").Value = 146.5
$ws.Range("J23").Value = 442.8
$ws.Range("K23").Value = 439.5
$ws.Range("L23").Value = 1328.4
$ws.Range("M23").Value = -204.5
$ws.Range("N23").Value = -1798.4
$ws.Range("H122").Value = 6523417
$ws.Range("J122").Value = 1352385.8
$ws.Range("L122").Value = 12171472.2
$ws.Range("N122").Value = -12176372.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 26521.625
$ws.Range("J80").Value = 34666
$ws.Range("L80").Value = 34666
$ws.Range("N80").Value = -36662
$ws.Range("H83").Value = 26521.625
$ws.Range("J83").Value = 34666
$ws.Range("L83").Value = 173330
$ws.Range("N83").Value = -183314
$ws.Range("H122").Value = 1786524.1
$ws.Range("I122").Value = 2262740.5
$ws.Range("J122").Value = 712.5
$ws.Range("K122").Value = 6788221.5
$ws.Range("L122").Value = 2137.5
$ws.Range("M122").Value = -6785771.5
$ws.Range("N122").Value = -7037.5
$ws.Range("H123").Value = 52929.668
$ws.Range("J123").Value = 52929.668
$ws.Range("L123").Value = 52929.668
$ws.Range("N123").Value = -57829.668
$ws.Range("H126").Value = 4885038
$ws.Range("I126").Value = 2572735.8
$ws.Range("J126").Value = 8931567
$ws.Range("K126").Value = 7718207.399999999
$ws.Range("L126").Value = 26794701
$ws.Range("M126").Value = -7715737.399999999
$ws.Range("N126").Value = -26799641
$ws.Range("H132").Value = 3974.1428
$ws.Range("I132").Value = 3974.1428
$ws.Range("K132").Value = 11922.4284
$ws.Range("M132").Value = -9392.428400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2359175.5
$ws.Range("J40").Value = 4909202.5
$ws.Range("L40").Value = 4909202.5
$ws.Range("N40").Value = -4909474.5
$ws.Range("H55").Value = 2419.6
$ws.Range("I55").Value = 2499
$ws.Range("K55").Value = 2499
$ws.Range("M55").Value = -2326
$ws.Range("H61").Value = 3241.6296
$ws.Range("I61").Value = 2414.1304
$ws.Range("J61").Value = 7999.75
$ws.Range("K61").Value = 2414.1304
$ws.Range("L61").Value = 7999.75
$ws.Range("M61").Value = -2212.1304
$ws.Range("N61").Value = -8403.75
$ws.Range("H68").Value = 4200
$ws.Range("I68").Value = 3000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2251
$ws.Range("H71").Value = 4200
$ws.Range("I71").Value = 3000
$ws.Range("K71").Value = 15000
$ws.Range("M71").Value = -11256
$ws.Range("H100").Value = 3249.8462
$ws.Range("I100").Value = 3481.375
$ws.Range("J100").Value = 2879.4
$ws.Range("K100").Value = 3481.375
$ws.Range("L100").Value = 2879.4
$ws.Range("M100").Value = -2940.375
$ws.Range("N100").Value = -3961.4
$ws.Range("H113").Value = 3241.6296
$ws.Range("I113").Value = 2414.1304
$ws.Range("J113").Value = 7999.75
$ws.Range("K113").Value = 2414.1304
$ws.Range("L113").Value = 7999.75
$ws.Range("M113").Value = -244.1304
$ws.Range("N113").Value = -12339.75
$ws.Range("H122").Value = 31189810
$ws.Range("I122").Value = 52077384
$ws.Range("K122").Value = 156232152
$ws.Range("M122").Value = -156229702
$ws.Range("H136").Value = 12666.872
$ws.Range("I136").Value = 10166.077
$ws.Range("K136").Value = 30498.231
$ws.Range("M136").Value = -27948.231

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 494.4
$ws.Range("I100").Value = 441.69232
$ws.Range("K100").Value = 883.38464
$ws.Range("M100").Value = -342.38464
$ws.Range("H122").Value = 486749.5
$ws.Range("I122").Value = 675751.1
$ws.Range("K122").Value = 2027253.3
$ws.Range("M122").Value = -2024803.3
$ws.Range("H126").Value = 6672361
$ws.Range("I126").Value = 6167.75
$ws.Range("J126").Value = 33337134
$ws.Range("K126").Value = 18503.25
$ws.Range("L126").Value = 100011402
$ws.Range("M126").Value = -16033.25
$ws.Range("N126").Value = -100016342
$ws.Range("H136").Value = 12566.171
$ws.Range("J136").Value = 40078.09
$ws.Range("L136").Value = 120234.27
$ws.Range("N136").Value = -125334.27
